$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"5.07435593108119E-05"
$ws.Range("E2").Value = [double]"5.07435593108119E-05"

# Row 3
$ws.Range("D3").Value = [double]"0.9999987510404045"
$ws.Range("E3").Value = [double]"0.9999987510404045"

# Row 4
$ws.Range("D4").Value = [double]"0.001257508691300802"
$ws.Range("E4").Value = [double]"0.001257508691300802"

# Row 5
$ws.Range("D5").Value = [double]"0.0003370935981485177"
$ws.Range("E5").Value = [double]"0.0003370935981485177"

# Row 6
$ws.Range("D6").Value = [double]"0.05082270881061299"
$ws.Range("E6").Value = [double]"0.05082270881061299"

# Row 7
$ws.Range("D7").Value = [double]"0.9979753499352461"
$ws.Range("E7").Value = [double]"0.002024650064753941"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"0.4994972685212419"
$ws.Range("E8").Value = [double]"0.5005027314787581"

# Row 9
$ws.Range("D9").Value = [double]"0.8662292693367433"
$ws.Range("E9").Value = [double]"0.1337707306632567"

# Row 10
$ws.Range("D10").Value = [double]"0.9907412615257087"
$ws.Range("E10").Value = [double]"0.009258738474291262"

# Row 11
$ws.Range("D11").Value = [double]"0.8451179720641857"
$ws.Range("E11").Value = [double]"0.1548820279358143"
$ws.Range("F11").Value = [double]"1.466437101364136"
$ws.Range("G11").Value = [double]"0.8"

# Row 12
$ws.Range("D12").Value = [double]"6.161428462720524E-08"
$ws.Range("E12").Value = [double]"6.161428462720524E-08"

# Row 13
$ws.Range("D13").Value = [double]"0.999999998842094"
$ws.Range("E13").Value = [double]"0.999999998842094"

# Row 14
$ws.Range("D14").Value = [double]"4.663320743194811E-05"
$ws.Range("E14").Value = [double]"4.663320743194811E-05"

# Row 15
$ws.Range("D15").Value = [double]"0.0002252178284814003"
$ws.Range("E15").Value = [double]"0.0002252178284814003"

# Row 16
$ws.Range("D16").Value = [double]"0.039419527049603"
$ws.Range("E16").Value = [double]"0.039419527049603"

# Row 17
$ws.Range("D17").Value = [double]"0.999996439971073"
$ws.Range("E17").Value = [double]"3.560028926985481E-06"

# Row 18
$ws.Range("D18").Value = [double]"0.5174829194037111"
$ws.Range("E18").Value = [double]"0.4825170805962889"

# Row 19
$ws.Range("D19").Value = [double]"0.9649667707068637"
$ws.Range("E19").Value = [double]"0.03503322929313635"

# Row 20
$ws.Range("D20").Value = [double]"0.9999993902871904"
$ws.Range("E20").Value = [double]"6.097128095738924E-07"

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = [double]"0.2261852967477547"
$ws.Range("E21").Value = [double]"0.7738147032522453"
$ws.Range("F21").Value = [double]"2.279798984527588"
$ws.Range("G21").Value = [double]"0.8"
